$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 308). Bump that date forward by one day (46074 -> 46075)
# for all of them, keeping the existing date formatting/style intact.
$ws.Range("C2:C308").Value = 46075
